$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "Total" column (K) with a SUM formula for each student row,
# filled down from K4 through K9 so Excel stores it as a shared formula.
$ws.Range("K4:K9").Formula = "=SUM(E4:J4)"
